$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links")

# Insert a new row at the top, shifting existing rows down by one.
$ws.Rows.Item(1).Insert()

# Populate the newly inserted row 1 with the new link entry.
$ws.Range("A1").Value = 30
$ws.Range("B1").Value = "Описание SQLite"
$ws.Range("C1").Value = "https://metanit.com/sql/sqlite/"
$ws.Range("D1").Value = "Описание SQLite"
